$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'-0.56%"
$ws.Range("G2").Value = "'13"
$ws.Range("D3").Value = "'35.74"
$ws.Range("E3").Value = "'-0.25%"
$ws.Range("G3").Value = "'13"
$ws.Range("D4").Value = "'5.038"
$ws.Range("E4").Value = "'-0.32%"
$ws.Range("G4").Value = "'13"
$ws.Range("D5").Value = "'0.07993"
$ws.Range("E5").Value = "'-1.36%"
$ws.Range("G5").Value = "'13"
$ws.Range("D6").Value = "'1.858"
$ws.Range("E6").Value = "'-5.00%"
$ws.Range("G6").Value = "'13"
$ws.Range("D7").Value = "'4.123"
$ws.Range("E7").Value = "'-1.12%"
$ws.Range("G7").Value = "'13"
$ws.Range("D8").Value = "'7.765"
$ws.Range("E8").Value = "'-0.08%"
$ws.Range("G8").Value = "'13"
$ws.Range("D9").Value = "'0.9238"
$ws.Range("E9").Value = "'-0.72%"
$ws.Range("G9").Value = "'13"
$ws.Range("D10").Value = "'0.1289"
$ws.Range("E10").Value = "'-5.45%"
$ws.Range("G10").Value = "'13"
$ws.Range("D11").Value = "'0.1893"
$ws.Range("E11").Value = "'-0.65%"
$ws.Range("G11").Value = "'13"
$ws.Range("D12").Value = "'0.09072"
$ws.Range("E12").Value = "'-2.41%"
$ws.Range("G12").Value = "'13"
$ws.Range("D13").Value = "'0.03407"
$ws.Range("E13").Value = "'-3.40%"
$ws.Range("G13").Value = "'13"
$ws.Range("D14").Value = "'0.09851"
$ws.Range("E14").Value = "'-0.14%"
$ws.Range("G14").Value = "'13"
$ws.Range("D15").Value = "'0.001403"
$ws.Range("E15").Value = "'-2.36%"
$ws.Range("G15").Value = "'13"
$ws.Range("D16").Value = "'0.006176"
$ws.Range("E16").Value = "'6.42%"
$ws.Range("G16").Value = "'13"
$ws.Range("D17").Value = "'3.849"
$ws.Range("E17").Value = "'8.01%"
$ws.Range("G17").Value = "'13"
$ws.Range("D18").Value = "'3.367"
$ws.Range("E18").Value = "'14.07%"
$ws.Range("G18").Value = "'13"
$ws.Range("E19").Value = "'-0.84%"
$ws.Range("G19").Value = "'13"
$ws.Range("D20").Value = "'0.1334"
$ws.Range("E20").Value = "'-0.80%"
$ws.Range("G20").Value = "'13"
$ws.Range("D21").Value = "'4.812"
$ws.Range("E21").Value = "'-1.30%"
$ws.Range("G21").Value = "'13"
$ws.Range("E22").Value = "'-11.63%"
$ws.Range("G22").Value = "'13"
$ws.Range("D23").Value = "'0.04408"
$ws.Range("E23").Value = "'0.14%"
$ws.Range("G23").Value = "'13"
$ws.Range("E24").Value = "'1.19%"
$ws.Range("G24").Value = "'13"
$ws.Range("D25").Value = "'0.004872"
$ws.Range("E25").Value = "'1.85%"
$ws.Range("G25").Value = "'13"
$ws.Range("G26").Value = "'13"
$ws.Range("E27").Value = "'-21.16%"
$ws.Range("G27").Value = "'13"
$ws.Range("E28").Value = "'42.35%"
$ws.Range("G28").Value = "'13"
$ws.Range("G29").Value = "'13"
$ws.Range("G30").Value = "'13"
$ws.Range("G31").Value = "'13"
$ws.Range("G32").Value = "'13"
$ws.Range("G33").Value = "'13"
$ws.Range("G34").Value = "'13"
$ws.Range("G35").Value = "'13"
$ws.Range("G36").Value = "'13"
$ws.Range("G37").Value = "'13"
$ws.Range("G38").Value = "'13"
$ws.Range("D39").Value = "'0.01934"
$ws.Range("E39").Value = "'-1.95%"
$ws.Range("G39").Value = "'13"
$ws.Range("D40").Value = "'0.05159"
$ws.Range("E40").Value = "'3.66%"
$ws.Range("G40").Value = "'13"
$ws.Range("D41").Value = "'0.007600"
$ws.Range("E41").Value = "'-0.33%"
$ws.Range("G41").Value = "'13"
$ws.Range("D42").Value = "'0.01016"
$ws.Range("E42").Value = "'-5.10%"
$ws.Range("G42").Value = "'13"
$ws.Range("E43").Value = "'-2.26%"
$ws.Range("G43").Value = "'13"
$ws.Range("D44").Value = "'0.002172"
$ws.Range("E44").Value = "'3.59%"
$ws.Range("G44").Value = "'13"
$ws.Range("D45").Value = "'0.009910"
$ws.Range("E45").Value = "'-8.16%"
$ws.Range("G45").Value = "'13"
$ws.Range("D46").Value = "'0.00006199"
$ws.Range("E46").Value = "'-2.66%"
$ws.Range("G46").Value = "'13"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("G47").Value = "'13"
$ws.Range("D48").Value = "'64.84"
$ws.Range("E48").Value = "'-0.18%"
$ws.Range("G48").Value = "'13"
$ws.Range("E49").Value = "'39.66%"
$ws.Range("G49").Value = "'13"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("G50").Value = "'13"
$ws.Range("E51").Value = "'0.21%"
$ws.Range("G51").Value = "'13"
